$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - PROJECTILE_MOSI (A,B,D,E first)
$ws.Range("A16").Value2 = "PROJECTILE_MOSI"
$ws.Range("B16").Value2 = "PA23"
$ws.Range("D16").Value2 = "SERCOM5 PAD0"
$ws.Range("E16").Value2 = "SPI ASYNC"

# Row 17 - PROJECTILE_SCK (A,B,D first)
$ws.Range("A17").Value2 = "PROJECTILE_SCK"
$ws.Range("B17").Value2 = "PA22"
$ws.Range("D17").Value2 = "SERCOM5 PAD1"
$ws.Range("E17").Value2 = "SPI ASYNC"

# Now fill in the Driver-signal column (F) for rows 16 and 17
$ws.Range("F16").Value2 = "MOSI"
$ws.Range("F17").Value2 = "SCK"

# Row 18 - PROJECTILE_MISO (A,B first)
$ws.Range("A18").Value2 = "PROJECTILE_MISO"
$ws.Range("B18").Value2 = "PA21"
$ws.Range("E18").Value2 = "SPI ASYNC"

# F18 written before D18
$ws.Range("F18").Value2 = "MISO"
$ws.Range("D18").Value2 = "SERCOM5 PAD3"

# Row 19 - PROJECTILE_FITTED_FB
$ws.Range("A19").Value2 = "PROJECTILE_FITTED_FB"
$ws.Range("B19").Value2 = "PA20"
$ws.Range("C19").Value2 = "GPIO in, pull down."

# Apply the same centered alignment style used by the rest of the data rows
# (cell style index 2 => horizontal="center")
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("D16").HorizontalAlignment = -4108
$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("F16").HorizontalAlignment = -4108

$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("D17").HorizontalAlignment = -4108
$ws.Range("E17").HorizontalAlignment = -4108
$ws.Range("F17").HorizontalAlignment = -4108

$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("B18").HorizontalAlignment = -4108
$ws.Range("D18").HorizontalAlignment = -4108
$ws.Range("E18").HorizontalAlignment = -4108
$ws.Range("F18").HorizontalAlignment = -4108

$ws.Range("A19").HorizontalAlignment = -4108
$ws.Range("B19").HorizontalAlignment = -4108
$ws.Range("C19").HorizontalAlignment = -4108

# Columns now need to be widened to fit the newly entered pin names / values
# (mirrors Excel's "best fit" column auto-sizing after the data entry above)
$ws.Columns("A").ColumnWidth = 29.5
$ws.Columns("B").ColumnWidth = 7.6666666666667
$ws.Columns("C").ColumnWidth = 20.6666666666667
$ws.Columns("D").ColumnWidth = 13.8333333333333
$ws.Columns("E").ColumnWidth = 12.3333333333333
$ws.Columns("F").ColumnWidth = 12

# Update selection to match the saved view state
$ws.Range("C20").Select()
